$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the serial date 45172 (2023-09-03) for every
# data row (2 through 387). Update it to 45175 (2023-09-06) for all rows,
# preserving existing formatting/style.
for ($r = 2; $r -le 387; $r++) {
    $ws.Cells.Item($r, 3).Value = 45175
}
